$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Release 0.10.0: duplicate the "0.9.5" testing-protocol sheet into a new
# "0.10.0" sheet (placed right after "0.9.5"), then update it with the
# results of the 0.10.0 test pass.
# ---------------------------------------------------------------------------

$source = $wb.Worksheets.Item("0.9.5")
$source.Copy($null, $source)

# The copy is inserted immediately after the source sheet.
$newSheet = $wb.Worksheets.Item($source.Index + 1)
$newSheet.Name = "0.10.0"
$newSheet.Activate()

# --- Header banner (row 1) ------------------------------------------------
$newSheet.Range("C1").Value = "Roman: 0.10.0 win 64-bit (2016-06-08), Windows 7 64-bit, Java jdk1.8.0_74 64-bit with ConfModel"

# --- Result-column highlight moves from row 11 to row 12 -------------------
# Row 11 used to carry the highlighted "Result" cell + a comment; in the new
# sheet the comment is gone and the cell keeps its font but loses its fill.
# Row 12 is the one that now gets the yellow highlight + a new comment, so
# copy the (still intact) row-11 formatting over to row 12 first.
$newSheet.Range("C11").Copy() | Out-Null
$newSheet.Range("C12").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = $false

# Now strip the fill from row 11's Result cell (keep font/border/alignment).
$newSheet.Range("C11").Interior.Pattern = -4142           # xlNone

# Clear the old row-11 comment, add the new row-12 comment.
$newSheet.Range("D11").Clear()
$newSheet.Range("D12").Value = "AlgorithmFamily is not shown in the Algorithm section after Algorithm is assigned to it. Only after restart."

# --- Row 9: new comment added ----------------------------------------------
$newSheet.Range("D9").Value = "Need to select another property field to make editor dirty"

# --- Row 15: comment removed, row height shrinks ---------------------------
$newSheet.Range("D15").Clear()
$newSheet.Rows(15).RowHeight = 55.5

# --- Row 16: comment removed, highlight removed -----------------------------
$newSheet.Range("D16").Clear()
$newSheet.Range("C17").Copy() | Out-Null                  # style 1 (no highlight)
$newSheet.Range("C16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 21: new comment + highlight + taller row ---------------------------
$newSheet.Range("C16").Copy() | Out-Null                  # style 22 (highlight) -- reuse a cell that already carries it
$newSheet.Range("C21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$newSheet.Range("D21").Value = "after changing constraint from pipeline editor full qualified names are shown in the constraint"
$newSheet.Rows(21).RowHeight = 30

# --- Row 26: comment removed, highlight removed, row height back to default -
$newSheet.Range("D26").Clear()
$newSheet.Range("C17").Copy() | Out-Null                  # style 1 (no highlight)
$newSheet.Range("C26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$newSheet.Rows(26).AutoFit()

# --- Row 30: comment removed ------------------------------------------------
$newSheet.Range("D30").Clear()

# --- View state: scroll the frozen pane down, select the last data row -----
$newSheet.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 2
$newSheet.Range("D31").Select()
